$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend date formatting (style) to the two newly appended rows (16 and 17)
$ws.Range("D16").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("D17").NumberFormat = $ws.Range("D15").NumberFormat

# Row 5
$ws.Range("D5").Value = 44413

# Row 6
$ws.Range("D6").Value = 44413
$ws.Range("M6").Value = 200

# Row 7
$ws.Range("D7").Value = 44412
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 2600
$ws.Range("O7").Value = 2700
$ws.Range("P7").Value = 2650
$ws.Range("Q7").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R7").Value = 'Provincia del Elquí'
$ws.Range("S7").Value = 2650
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44412
$ws.Range("L8").Value = 'Segunda'
$ws.Range("N8").Value = 2200
$ws.Range("O8").Value = 2300
$ws.Range("P8").Value = 2250
$ws.Range("Q8").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("R8").Value = 'Provincia del Elquí'
$ws.Range("S8").Value = 2250
$ws.Range("T8").Value = 1

# Row 9
$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 13000
$ws.Range("O9").Value = 13500
$ws.Range("P9").Value = 13250
$ws.Range("S9").Value = 1656

# Row 10
$ws.Range("D10").Value = 44161
$ws.Range("L10").Value = 'Primera'
$ws.Range("N10").Value = 11000
$ws.Range("O10").Value = 11500
$ws.Range("P10").Value = 11250
$ws.Range("S10").Value = 1406

# Row 11
$ws.Range("D11").Value = 44161
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 9000
$ws.Range("O11").Value = 9500
$ws.Range("P11").Value = 9250
$ws.Range("S11").Value = 1156

# Row 12
$ws.Range("L12").Value = 'Especial'
$ws.Range("N12").Value = 12500
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 12750
$ws.Range("S12").Value = 1594

# Row 13
$ws.Range("D13").Value = 44160
$ws.Range("L13").Value = 'Primera'
$ws.Range("N13").Value = 10500
$ws.Range("O13").Value = 11000
$ws.Range("P13").Value = 10750
$ws.Range("S13").Value = 1344

# Row 14
$ws.Range("D14").Value = 44160
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 240
$ws.Range("N14").Value = 8500
$ws.Range("O14").Value = 9000
$ws.Range("P14").Value = 8750
$ws.Range("S14").Value = 1094

# Row 15
$ws.Range("L15").Value = 'Especial'
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 13000
$ws.Range("O15").Value = 13500
$ws.Range("P15").Value = 13250
$ws.Range("S15").Value = 1656

# Row 16
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C16").Value = 'Coquimbo'
$ws.Range("D16").Value = 44168
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 'Fruta'
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = 'Otros'
$ws.Range("I16").Value = 100107002
$ws.Range("J16").Value = 'Chirimoya'
$ws.Range("K16").Value = 'Cultivar IV Región'
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 11000
$ws.Range("O16").Value = 11500
$ws.Range("P16").Value = 11250
$ws.Range("Q16").Value = '$/bandeja 8 kilos'
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1406
$ws.Range("T16").Value = 8

# Row 17
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C17").Value = 'Coquimbo'
$ws.Range("D17").Value = 44168
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 'Fruta'
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = 'Otros'
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = 'Chirimoya'
$ws.Range("K17").Value = 'Cultivar IV Región'
$ws.Range("L17").Value = 'Segunda'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 8500
$ws.Range("O17").Value = 9000
$ws.Range("P17").Value = 8750
$ws.Range("Q17").Value = '$/bandeja 8 kilos'
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 1094
$ws.Range("T17").Value = 8

